# Integrating Derek's viz notebooks into my branch

$wb = $excel.ActiveWorkbook

# --- "owners" sheet: update scenario values & damage states ---
$owners = $wb.Worksheets.Item("owners")

$owners.Range("C2").Value = 5000
$owners.Range("J2").Value = 1100

$owners.Range("H3").Value = 4
$owners.Range("M3").Value = "Moderate"

$owners.Range("C4").Value = 2500
$owners.Range("H4").Value = 2
$owners.Range("J4").Value = 1200
$owners.Range("M4").Value = "Extensive"

$owners.Range("H5").Value = 3
$owners.Range("M5").Value = "Extensive"

# Leave the "owners" sheet selection on M5 (matches the saved selection state)
$owners.Range("M5").Select()

# --- "human_capital" sheet: updated workforce quantities ---
$hc = $wb.Worksheets.Item("human_capital")

$hc.Range("B2").Value = 2
$hc.Range("B3").Value = 2
$hc.Range("B4").Value = 1
$hc.Range("B8").Value = 2

# Make human_capital the active sheet/tab with B9 selected
$hc.Activate()
$hc.Range("B9").Select()
